$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values per corrected
# Diebold-Mariano computations.

$ws.Range("C2").Value = 1.089444305578694
$ws.Range("D2").Value = 0.2877481520155922

$ws.Range("C3").Value = 1.032393118162679
$ws.Range("D3").Value = 0.3131044770999618

$ws.Range("C4").Value = 0.1637249140059307
$ws.Range("D4").Value = 0.8714424993452705

$ws.Range("C5").Value = -0.2255424715097639
$ws.Range("D5").Value = 0.8236393153408452

$ws.Range("C6").Value = 0.02928714908817166
$ws.Range("D6").Value = 0.9768996378920072

$ws.Range("C7").Value = -0.9078757486611765
$ws.Range("D7").Value = 0.3737819452636848

$ws.Range("C8").Value = -0.9392823303279144
$ws.Range("D8").Value = 0.3577848825201915

$ws.Range("C9").Value = -1.030485535004753
$ws.Range("D9").Value = 0.313978718674363

$ws.Range("C10").Value = -1.003494111983461
$ws.Range("D10").Value = 0.3265329166393238

$ws.Range("C11").Value = -0.3937967733850236
$ws.Range("D11").Value = 0.6975219132791275
